# "fixes: owner changes in attribute"
#
# The "Owner" column (column D) on the request-side sheets (publish, revise,
# query) was uniformly (and incorrectly) set to "Publisher" for every
# attribute row. This fixes it so those rows are owned by the "Provider"
# instead. The callback/response sheets (on_publish, on_revise, on_query)
# already had the correct, mixed Publisher/Provider values per row and are
# left unchanged.

$wb = $excel.ActiveWorkbook

$sheetNames = @("publish", "revise", "query")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find last used row in column A (Path column) to know how far to scan.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        $val = $cell.Value2

        if ($val -eq "Publisher") {
            $cell.Value = "Provider"
        }
    }
}
